# Auto-generated Excel COM-interop script
# Applies the Sheets/Twintania_Profits.xlsx cell-value updates described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 882.55554
$ws.Range("I6").Value = 268
$ws.Range("K6").Value = 804
$ws.Range("M6").Value = -692
$ws.Range("H40").Value = 1798.909
$ws.Range("I40").Value = 1798.909
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1798.909
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1623.909
$ws.Range("H41").Value = 721.1429000000001
$ws.Range("I41").Value = 612.25
$ws.Range("J41").Value = 866.3333
$ws.Range("K41").Value = 612.25
$ws.Range("L41").Value = 866.3333
$ws.Range("M41").Value = -172.25
$ws.Range("N41").Value = -1746.3333
$ws.Range("H43").Value = 3484
$ws.Range("I43").Value = 3381.3333
$ws.Range("K43").Value = 3381.3333
$ws.Range("M43").Value = -3312.3333
$ws.Range("H86").Value = 6542.909
$ws.Range("I86").Value = 2272.75
$ws.Range("K86").Value = 2272.75
$ws.Range("M86").Value = -1149.75
$ws.Range("H89").Value = 6542.909
$ws.Range("I89").Value = 2272.75
$ws.Range("K89").Value = 11363.75
$ws.Range("M89").Value = -5747.75
$ws.Range("H132").Value = 2280.8
$ws.Range("I132").Value = 2226
$ws.Range("K132").Value = 6678
$ws.Range("M132").Value = -4148
$ws.Range("H137").Value = 12856.974
$ws.Range("I137").Value = 5472.5264
$ws.Range("J137").Value = 20241.422
$ws.Range("K137").Value = 16417.5792
$ws.Range("L137").Value = 60724.266
$ws.Range("M137").Value = -13867.5792
$ws.Range("N137").Value = -65824.266
$ws.Range("H138").Value = 3077.1
$ws.Range("I138").Value = 2688.4
$ws.Range("J138").Value = 4243.2
$ws.Range("K138").Value = 8065.200000000001
$ws.Range("L138").Value = 12729.6
$ws.Range("M138").Value = -2925.200000000001
$ws.Range("N138").Value = -23009.6
$ws.Range("H141").Value = 3279.2
$ws.Range("I141").Value = 3156.2856
$ws.Range("K141").Value = 9468.856800000001
$ws.Range("M141").Value = -4288.856800000001
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4078.1082
$ws.Range("I2").Value = 3253.3
$ws.Range("K2").Value = 3253.3
$ws.Range("M2").Value = -3140.3
$ws.Range("H32").Value = 2269.191
$ws.Range("I32").Value = 1544.65
$ws.Range("K32").Value = 1544.65
$ws.Range("M32").Value = -1257.65
$ws.Range("H61").Value = 7705.4116
$ws.Range("I61").Value = 5586.6
$ws.Range("K61").Value = 5586.6
$ws.Range("M61").Value = -5374.6
$ws.Range("H74").Value = 2360.3333
$ws.Range("I74").Value = 1678.2667
$ws.Range("K74").Value = 1678.2667
$ws.Range("M74").Value = -804.2666999999999
$ws.Range("H77").Value = 2360.3333
$ws.Range("I77").Value = 1678.2667
$ws.Range("K77").Value = 8391.333499999999
$ws.Range("M77").Value = -4023.333499999999
$ws.Range("H110").Value = 1048.25
$ws.Range("I110").Value = 740
$ws.Range("K110").Value = 740
$ws.Range("M110").Value = 1305
$ws.Range("H114").Value = 76000
$ws.Range("J114").Value = 76000
$ws.Range("L114").Value = 76000
$ws.Range("N114").Value = -84678
$ws.Range("H116").Value = 4078.1082
$ws.Range("I116").Value = 3253.3
$ws.Range("K116").Value = 3253.3
$ws.Range("M116").Value = -959.3000000000002
$ws.Range("H122").Value = 740
$ws.Range("I122").Value = 740
$ws.Range("K122").Value = 2220
$ws.Range("M122").Value = 230
$ws.Range("H136").Value = 7705.4116
$ws.Range("I136").Value = 5586.6
$ws.Range("K136").Value = 16759.8
$ws.Range("M136").Value = -14209.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4078.1082
$ws.Range("I3").Value = 3253.3
$ws.Range("K3").Value = 3253.3
$ws.Range("M3").Value = -3139.3
$ws.Range("H134").Value = 7991.5938
$ws.Range("I134").Value = 4406
$ws.Range("K134").Value = 13218
$ws.Range("M134").Value = -10683

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1427.3636
$ws.Range("I19").Value = 1427.3636
$ws.Range("K19").Value = 1427.3636
$ws.Range("M19").Value = -1257.3636
$ws.Range("H24").Value = 1427.3636
$ws.Range("I24").Value = 1427.3636
$ws.Range("K24").Value = 1427.3636
$ws.Range("M24").Value = -1257.3636
$ws.Range("H25").Value = 12686
$ws.Range("I25").Value = 12686
$ws.Range("K25").Value = 12686
$ws.Range("M25").Value = -12512
$ws.Range("H62").Value = 335412.12
$ws.Range("I62").Value = 169082.83
$ws.Range("J62").Value = 668070.7
$ws.Range("K62").Value = 169082.83
$ws.Range("L62").Value = 668070.7
$ws.Range("M62").Value = -168458.83
$ws.Range("N62").Value = -669318.7
$ws.Range("H65").Value = 335412.12
$ws.Range("I65").Value = 169082.83
$ws.Range("J65").Value = 668070.7
$ws.Range("K65").Value = 845414.1499999999
$ws.Range("L65").Value = 3340353.5
$ws.Range("M65").Value = -842294.1499999999
$ws.Range("N65").Value = -3346593.5
$ws.Range("H99").Value = 12700.186
$ws.Range("I99").Value = 12097.5
$ws.Range("K99").Value = 12097.5
$ws.Range("M99").Value = -10599.5
$ws.Range("H126").Value = 12700.186
$ws.Range("I126").Value = 12097.5
$ws.Range("K126").Value = 36292.5
$ws.Range("M126").Value = -33822.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5500
$ws.Range("I3").Value = 5500
$ws.Range("K3").Value = 16500
$ws.Range("M3").Value = -16388
$ws.Range("H34").Value = 2709.68
$ws.Range("I34").Value = 114.35294
$ws.Range("J34").Value = 8224.75
$ws.Range("K34").Value = 343.05882
$ws.Range("L34").Value = 24674.25
$ws.Range("M34").Value = -259.05882
$ws.Range("N34").Value = -24842.25
$ws.Range("H55").Value = 2071.4546
$ws.Range("J55").Value = 2928.2856
$ws.Range("L55").Value = 8784.856800000001
$ws.Range("N55").Value = -9138.856800000001
$ws.Range("H81").Value = 5586.8887
$ws.Range("I81").Value = 5956.6
$ws.Range("J81").Value = 5124.75
$ws.Range("K81").Value = 17869.8
$ws.Range("L81").Value = 15374.25
$ws.Range("M81").Value = -16746.8
$ws.Range("N81").Value = -17620.25
$ws.Range("H84").Value = 5586.8887
$ws.Range("I84").Value = 5956.6
$ws.Range("J84").Value = 5124.75
$ws.Range("K84").Value = 53609.4
$ws.Range("L84").Value = 46122.75
$ws.Range("M84").Value = -47993.4
$ws.Range("N84").Value = -57354.75
$ws.Range("H113").Value = 1149.8
$ws.Range("J113").Value = 728.5
$ws.Range("L113").Value = 2185.5
$ws.Range("N113").Value = -6525.5
$ws.Range("H134").Value = 895.087
$ws.Range("I134").Value = 895.087
$ws.Range("K134").Value = 2685.261
$ws.Range("M134").Value = 2384.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8500.25
$ws.Range("I10").Value = 2501
$ws.Range("K10").Value = 2501
$ws.Range("M10").Value = -2332
$ws.Range("H122").Value = 1793
$ws.Range("I122").Value = 1793
$ws.Range("K122").Value = 5379
$ws.Range("M122").Value = -2929
$ws.Range("H126").Value = 5336.6875
$ws.Range("I126").Value = 3854.5
$ws.Range("K126").Value = 11563.5
$ws.Range("M126").Value = -9093.5
$ws.Range("H132").Value = 3200.16
$ws.Range("I132").Value = 3105.4736
$ws.Range("K132").Value = 9316.4208
$ws.Range("M132").Value = -6786.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 8317.809999999999
$ws.Range("I25").Value = 8104.625
$ws.Range("J25").Value = 9000
$ws.Range("K25").Value = 8104.625
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = -7874.625
$ws.Range("N25").Value = -9460
$ws.Range("H61").Value = 1481.9584
$ws.Range("I61").Value = 1471.4375
$ws.Range("K61").Value = 1471.4375
$ws.Range("M61").Value = -1269.4375
$ws.Range("H113").Value = 1481.9584
$ws.Range("I113").Value = 1471.4375
$ws.Range("K113").Value = 1471.4375
$ws.Range("M113").Value = 698.5625
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H122").Value = 4258
$ws.Range("I122").Value = 3099.6667
$ws.Range("J122").Value = 5416.3335
$ws.Range("K122").Value = 9299.000100000001
$ws.Range("L122").Value = 16249.0005
$ws.Range("M122").Value = -6849.000100000001
$ws.Range("N122").Value = -21149.0005
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 55596.25
$ws.Range("J109").Value = 55596.25
$ws.Range("L109").Value = 55596.25
$ws.Range("N109").Value = -58370.25
$ws.Range("H122").Value = 6303.3335
$ws.Range("I122").Value = 3545.9167
$ws.Range("K122").Value = 10637.7501
$ws.Range("M122").Value = -8187.750100000001
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H132").Value = 22313.488
$ws.Range("I132").Value = 16568.6
$ws.Range("J132").Value = 29494.6
$ws.Range("K132").Value = 49705.8
$ws.Range("L132").Value = 88483.79999999999
$ws.Range("M132").Value = -47175.8
$ws.Range("N132").Value = -93543.79999999999
